# Castor_validation.xlsx — "upgrade dfa summary function to calculate visit
# times diff": the visit-date summary now computes additional SV_n_Date
# entries for several subjects (extending each subject's visit-date series
# with the newly-calculated next visit date(s)). No existing values change;
# we only fill in previously-empty date cells (and one status cell) further
# along each subject's row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# SubjectID 206-011 (row 12): add SV_3..SV_6 dates
$ws.Range("G12").Value = "21-08-2023"
$ws.Range("H12").Value = "28-08-2023"
$ws.Range("I12").Value = "04-09-2023"
$ws.Range("J12").Value = "14-09-2023"

# SubjectID 207-011 (row 23): add SV_5..SV_6 dates
$ws.Range("I23").Value = "05-09-2023"
$ws.Range("J23").Value = "12-09-2023"

# SubjectID 207-013 (row 25): add SV_2..SV_3 dates
$ws.Range("F25").Value = "05-09-2023"
$ws.Range("G25").Value = "12-09-2023"

# SubjectID 208-016 (row 41): add SV_10 date
$ws.Range("N41").Value = "06-09-2023"

# SubjectID 208-018 (row 43): add SV_8 date
$ws.Range("L43").Value = "06-09-2023"

# SubjectID 209-011 (row 60): add SV_8..SV_10 dates
$ws.Range("L60").Value = "26-06-2023"
$ws.Range("M60").Value = "05-07-2023"
$ws.Range("N60").Value = "10-07-2023"

# SubjectID 209-014 (row 63): add SV_4..SV_5 dates
$ws.Range("H63").Value = "05-09-2023"
$ws.Range("I63").Value = "12-09-2023"

# SubjectID 209-015 (row 64): add SV_4..SV_5 dates
$ws.Range("H64").Value = "06-09-2023"
$ws.Range("I64").Value = "13-09-2023"

# SubjectID 209-019 (row 68): add SV_1..SV_3 dates
$ws.Range("F68").Value = "29-08-2023"
$ws.Range("G68").Value = "05-09-2023"
$ws.Range("H68").Value = "12-09-2023"

# SubjectID 210-003 (row 71): add SV_6 date
$ws.Range("K71").Value = "01-06-2023"

# SubjectID 210-004 (row 72): add SV_5..SV_6 dates
$ws.Range("I72").Value = "25-05-2023"
$ws.Range("J72").Value = "01-06-2023"

# SubjectID 210-005 (row 73): status now known, plus SV_1..SV_2 dates
$ws.Range("C73").Value = "In progress"
$ws.Range("E73").Value = "27-07-2023"
$ws.Range("F73").Value = "02-08-2023"

# SubjectID 211-004 (row 77): add SV_2 date
$ws.Range("F77").Value = "11-08-2023"

# SubjectID 211-005 (row 78): add SV_4..SV_5 dates
$ws.Range("H78").Value = "08-09-2023"
$ws.Range("I78").Value = "15-09-2023"

# SubjectID 211-006 (row 79): add SV_3..SV_5 dates
$ws.Range("G79").Value = "30-08-2023"
$ws.Range("H79").Value = "05-09-2023"
$ws.Range("I79").Value = "12-09-2023"
